$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSVData")

# Delete row 1 entirely, shifting all data rows up by one (oldest row drops off,
# formatting travels with the shifted rows, matching how Excel's Delete Row works).
$ws.Rows("1").Delete()

# Refresh the four newest transaction dates (formerly rows 2-5, now rows 1-4)
# with newer dates from the bank feed.
$ws.Range("A1").Value = 44540
$ws.Range("A2").Value = 44527
$ws.Range("A3").Value = 44491
$ws.Range("A4").Value = 44467

# Update the selected cell to reflect the new cursor position after the edit.
$ws.Range("A4").Select()
